# Add a new data row (row 3) to the "Planilha de Usuários" sheet, matching
# the formatting already used by the existing data row (row 2), and widen
# a couple of columns so the new, longer values fit comfortably.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start the new row off with the same formatting (font/alignment/etc.) as
# the existing data row, then fill in the values.
$ws.Range("A2:H2").Copy()
$ws.Range("A3:H3").PasteSpecial(-4122)

$ws.Range("A3").Value = 2
$ws.Range("B3").Value = "Marineia"
$ws.Range("C3").Value = "Almeida"
$ws.Range("D3").Value = "marineia123@gmail.com"
$ws.Range("E3").Value = 19
$ws.Range("F3").Value = "F"

# The CPF needs to stay a text value (it can carry leading zeros), so build
# it as a text formula and then flatten it back down to a plain value -
# that keeps the cell's existing (General) number format/style intact
# instead of Excel silently converting the digits to a number.
$ws.Range("G3").Formula = "=""12345678917"""
$ws.Range("G3").Copy()
$ws.Range("G3").PasteSpecial(-4163)

$ws.Range("H3").Value = "Aa123456789*"

# --- Column widths ---------------------------------------------------------
# "Nome" and "Email" need a bit more room for the new entries. ColumnWidth
# is expressed in characters, and Excel pads it by 5/6 of a character when
# storing the sheet's raw column width, so back that padding out up front
# to land on exact target widths of 10 and 23 characters.
$ws.Columns.Item(2).ColumnWidth = 10 - 5/6
$ws.Columns.Item(4).ColumnWidth = 23 - 5/6
